$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column O holds "timestamp"; find last used row dynamically (xlUp = -4162)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 15).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 15).Value = "2023-02-21 12:56:52"
}
